$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B - shifts existing B..S to C..T
$ws.Range("B1").EntireColumn.Insert()

# New column B width (matches column A's rendered width, not best-fit)
$ws.Range("B1:B1").EntireColumn.ColumnWidth = 8.43

# Row 12: new scenario name/radio-button values (columns E/F after the shift)
# (set first so new shared strings are appended in the same order as the authored file)
$ws.Range("E12").Value = "NewImportLogic_2 - Test_Automation_2"
$ws.Range("F12").Value = "NewImportLogic_2 - Test_Automation_2_radio_button"

# New Usertype column header + values
$ws.Range("B1").Value = "Usertype"
$ws.Range("B2").Value = "Admin User"
$ws.Range("B7").Value = "Staff User"
$ws.Range("B12").Value = "Client User"

# Update selection to match authored state
$ws.Range("B12").Select()
